$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0

$ws.Range("H94").Value = 4309.0527
$ws.Range("I94").Value = 2374.7693
$ws.Range("J94").Value = 8500
$ws.Range("K94").Value = 2374.7693
$ws.Range("L94").Value = 8500
$ws.Range("M94").Value = -1923.7693
$ws.Range("N94").Value = -9402

$ws.Range("H101").Value = 83333544
$ws.Range("I101").Value = 83333544
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 250000632
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -249999010

$ws.Range("H103").Value = 441.42856
$ws.Range("I103").Value = 364.66666
$ws.Range("J103").Value = 499
$ws.Range("K103").Value = 1093.99998
$ws.Range("L103").Value = 1497
$ws.Range("M103").Value = -507.9999800000001
$ws.Range("N103").Value = -2669

$ws.Range("H116").Value = 2250
$ws.Range("I116").Value = 2250
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2250
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1192

$ws.Range("H129").Value = 1293.88
$ws.Range("I129").Value = 870
$ws.Range("J129").Value = 1311.5416
$ws.Range("K129").Value = 2610
$ws.Range("L129").Value = 3934.6248
$ws.Range("M129").Value = 2390
$ws.Range("N129").Value = -13934.6248

$ws.Range("H132").Value = 1468.56
$ws.Range("I132").Value = 1427.102
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 4281.306
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -1751.306
$ws.Range("N132").Value = -15560

$ws.Range("H138").Value = 2382.19
$ws.Range("I138").Value = 1187.6923
$ws.Range("J138").Value = 4600.543
$ws.Range("K138").Value = 3563.0769
$ws.Range("L138").Value = 13801.629
$ws.Range("M138").Value = 1576.9231
$ws.Range("N138").Value = -24081.629

$ws.Range("H141").Value = 5357.636
$ws.Range("I141").Value = 1182.6888
$ws.Range("J141").Value = 24144.9
$ws.Range("K141").Value = 3548.0664
$ws.Range("L141").Value = 72434.70000000001
$ws.Range("M141").Value = 1631.9336
$ws.Range("N141").Value = -82794.70000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N6").Value = ""
$ws.Range("H6").Value = 38251.5
$ws.Range("I6").Value = 38251.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 38251.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -38078.5

$ws.Range("H32").Value = 5797.472
$ws.Range("I32").Value = 4094.3823
$ws.Range("J32").Value = 34750
$ws.Range("K32").Value = 4094.3823
$ws.Range("L32").Value = 34750
$ws.Range("M32").Value = -3807.3823
$ws.Range("N32").Value = -35324

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 46268
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 46268
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 46268
$ws.Range("N132").Value = -56388

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N15").Value = ""
$ws.Range("H15").Value = 9000
$ws.Range("I15").Value = 9000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -8830

$ws.Range("H16").Value = 3007
$ws.Range("I16").Value = 3269.625
$ws.Range("J16").Value = 1956.5
$ws.Range("K16").Value = 3269.625
$ws.Range("L16").Value = 1956.5
$ws.Range("M16").Value = -2982.625
$ws.Range("N16").Value = -2530.5

$ws.Range("H31").Value = 2302.1487
$ws.Range("I31").Value = 1368.5
$ws.Range("J31").Value = 3400.5588
$ws.Range("K31").Value = 1368.5
$ws.Range("L31").Value = 3400.5588
$ws.Range("M31").Value = -1073.5
$ws.Range("N31").Value = -3990.5588

$ws.Range("H34").Value = 2302.1487
$ws.Range("I34").Value = 1368.5
$ws.Range("J34").Value = 3400.5588
$ws.Range("K34").Value = 1368.5
$ws.Range("L34").Value = 3400.5588
$ws.Range("M34").Value = -1166.5
$ws.Range("N34").Value = -3804.5588

$ws.Range("M36").Value = ""
$ws.Range("H36").Value = 50000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 50000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 50000
$ws.Range("N36").Value = -50776

$ws.Range("M40").Value = ""
$ws.Range("H40").Value = 50000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 50000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 50000
$ws.Range("N40").Value = -50320

$ws.Range("H42").Value = 50000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 50000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -51186

$ws.Range("H58").Value = 1395.6719
$ws.Range("I58").Value = 1071.88
$ws.Range("J58").Value = 2552.0715
$ws.Range("K58").Value = 1071.88
$ws.Range("L58").Value = 2552.0715
$ws.Range("M58").Value = -868.8800000000001
$ws.Range("N58").Value = -2958.0715

$ws.Range("H113").Value = 3007
$ws.Range("I113").Value = 3269.625
$ws.Range("J113").Value = 1956.5
$ws.Range("K113").Value = 3269.625
$ws.Range("L113").Value = 1956.5
$ws.Range("M113").Value = -1099.625
$ws.Range("N113").Value = -6296.5

$ws.Range("H132").Value = 2068.3103
$ws.Range("I132").Value = 1052.4117
$ws.Range("J132").Value = 3507.5
$ws.Range("K132").Value = 3157.2351
$ws.Range("L132").Value = 10522.5
$ws.Range("M132").Value = -627.2351000000003
$ws.Range("N132").Value = -15582.5

$ws.Range("H136").Value = 1395.6719
$ws.Range("I136").Value = 1071.88
$ws.Range("J136").Value = 2552.0715
$ws.Range("K136").Value = 3215.64
$ws.Range("L136").Value = 7656.2145
$ws.Range("M136").Value = -665.6400000000003
$ws.Range("N136").Value = -12756.2145

$ws.Range("H140").Value = 54834.363
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54834.363
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54834.363
$ws.Range("N140").Value = -65194.363

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 712131.75
$ws.Range("I107").Value = 1619.5
$ws.Range("J107").Value = 1228867.9
$ws.Range("K107").Value = 4858.5
$ws.Range("L107").Value = 3686603.7
$ws.Range("M107").Value = -2938.5
$ws.Range("N107").Value = -3690443.7

$ws.Range("H122").Value = 1121.1154
$ws.Range("I122").Value = 430.76923
$ws.Range("J122").Value = 1811.4615
$ws.Range("K122").Value = 3876.92307
$ws.Range("L122").Value = 16303.1535
$ws.Range("M122").Value = -1426.92307
$ws.Range("N122").Value = -21203.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M17").Value = ""
$ws.Range("H17").Value = 12200
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 12200
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 12200
$ws.Range("N17").Value = -12536

$ws.Range("H132").Value = 3235.3333
$ws.Range("I132").Value = 3832
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 11496
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -8966
$ws.Range("N132").Value = -12260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 7073.4707
$ws.Range("I93").Value = 9270.416999999999
$ws.Range("J93").Value = 1800.8
$ws.Range("K93").Value = 9270.416999999999
$ws.Range("L93").Value = 1800.8
$ws.Range("M93").Value = -8022.416999999999

$ws.Range("H132").Value = 8873.325000000001
$ws.Range("I132").Value = 11105.565
$ws.Range("J132").Value = 6306.25
$ws.Range("K132").Value = 33316.695
$ws.Range("L132").Value = 18918.75
$ws.Range("M132").Value = -30786.695
$ws.Range("N132").Value = -23978.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 70007
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 70007
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 70007
$ws.Range("N15").Value = -70583

$ws.Range("H18").Value = 70007
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 70007
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 70007
$ws.Range("N18").Value = -70353

$ws.Range("H20").Value = 17003.666

$ws.Range("H107").Value = 2489.5
$ws.Range("I107").Value = 3168.5715
$ws.Range("J107").Value = 1810.4286
$ws.Range("K107").Value = 9505.7145
$ws.Range("L107").Value = 5431.2858
$ws.Range("M107").Value = -7585.7145
$ws.Range("N107").Value = -9271.2858

$ws.Range("H132").Value = 3029.8
$ws.Range("I132").Value = 3103.3572
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 9310.071599999999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -6780.071599999999
$ws.Range("N132").Value = -11060
